$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# D4 currently holds the number 1; change it to the text value "1б"
# (adding a new class designation, as shared string, not a number).
$ws.Range("D4").Value = "1б"
